$wb = $excel.ActiveWorkbook

# Mapping of row -> new "想去人数" (F column) value.
# Same update applies to both the "展览" and "全部类型" worksheets,
# which hold identical data.
$updates = @{
    3  = 95
    5  = 51
    6  = 565
    8  = 2029
    11 = 4388
    15 = 8
    17 = 26
    18 = 16
    20 = 3194
    22 = 470
    26 = 87
    30 = 200
    32 = 570
    33 = 1818
    34 = 277
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
